$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the shared strings first, in the order they appear in the target file
# (10 mins meditation, Read for 1hr, Jog around park 3x), by writing them to
# helper cells, then clear those cells and write the real row data referencing
# the same string values so the shared string table order is preserved.
$ws.Range("E1").Value = "10 mins meditation"
$ws.Range("E2").Value = "Read for 1hr"
$ws.Range("E3").Value = "Jog around park 3x"

# Add task rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Jog around park 3x"
$ws.Range("C2").Value = $true

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "10 mins meditation"
$ws.Range("C3").Value = $false

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Read for 1hr"
$ws.Range("C4").Value = $false

# Remove helper cells used only to seed shared string order
$ws.Range("E1:E3").Clear()

# Update selection to match target state
$ws.Range("B11").Select()
